# Update New Orleans shard workbook:
#  1. hotel_info gains a new "State" column (value "Louisiana") inserted
#     between "Hotel_Name" and "City".
#  2. The sheet tabs are reordered so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new State column into hotel_info -----------------------
$hotelWs = $wb.Worksheets.Item("hotel_info")
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Move review_info before hotel_info ---------------------------------
$reviewWs = $wb.Worksheets.Item("review_info")
$reviewWs.Move($wb.Worksheets.Item(1))
